$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": update handoff/handback datetimes (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-17 14:59:27"
$wsZhCn.Range("G2").Value = "2016-01-17 15:00:24"

# Sheet "de-de": update handoff/handback datetimes (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-17 14:59:38"
$wsDeDe.Range("G2").Value = "2016-01-17 15:00:43"
